$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
SELECT
    COUNT(DISTINCT std.dbgap_accession) AS "Studies",
    COUNT(DISTINCT prt.participant_id) AS "Participants",
    COUNT(DISTINCT smp.sample_id) AS "Samples",
    COALESCE(
        COUNT(DISTINCT sequ.id) + 
        COUNT(DISTINCT methyl.id) + 
        COUNT(DISTINCT patho.id) + 
        COUNT(DISTINCT cmf.id),
        0
    ) AS "Files"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
LEFT JOIN 
    df_diagnosis diag ON prt.id = diag."participant.id"     
LEFT JOIN 
    df_clinical_measure_file cmf ON prt.id = cmf."participant.id"
LEFT JOIN 
    df_sequencing_file sequ ON smp.id = sequ."sample.id"
LEFT JOIN 
    df_methylation_array_file methyl ON smp.id = methyl."sample.id"
LEFT JOIN 
    df_pathology_file patho ON smp.id = patho."sample.id"
WHERE 
    std.dbgap_accession = 'phs002790';
'@

$ws.Range("C2").Value = $newQuery

$ws.Activate()
$ws.Range("C3").Select()
$excel.ActiveWindow.Zoom = 70
